# "SwapJD >> MPShufleJD" -- rename the SwapJD sheet to MPShuffleJD and make
# it the active tab (it was previously Top1Absent).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SwapJD")
$ws.Name = "MPShuffleJD"

# The workbook carries a flat, hidden "_xlchart.v1.N" defined-name table that
# backs every box-and-whisker chart's data ranges. Before the edit it held
# three duplicated blocks (MPublished x2, Top1Absent x3) in addition to the
# single blocks used by the SwapJD and Top1Absent charts. Renaming SwapJD
# collapses/renumbers that table (Excel drops the unused duplicates and
# reassigns contiguous indices): the slot that used to be a spare MPublished
# duplicate (.20-.29) becomes the renamed sheet's block, and the slot that
# used to belong to SwapJD (.30-.39) is reused for Top1Absent, so the trailing
# duplicate Top1Absent blocks (.40-.69) disappear entirely.

# .20-.29: was a duplicate of MPublished!$B$1:$F$34 -> now MPShuffleJD's block
$wb.Names.Item("_xlchart.v1.20").RefersTo = "=MPShuffleJD!`$B`$1"
$wb.Names.Item("_xlchart.v1.21").RefersTo = "=MPShuffleJD!`$B`$2:`$B`$34"
$wb.Names.Item("_xlchart.v1.22").RefersTo = "=MPShuffleJD!`$C`$1"
$wb.Names.Item("_xlchart.v1.23").RefersTo = "=MPShuffleJD!`$C`$2:`$C`$34"
$wb.Names.Item("_xlchart.v1.24").RefersTo = "=MPShuffleJD!`$D`$1"
$wb.Names.Item("_xlchart.v1.25").RefersTo = "=MPShuffleJD!`$D`$2:`$D`$34"
$wb.Names.Item("_xlchart.v1.26").RefersTo = "=MPShuffleJD!`$E`$1"
$wb.Names.Item("_xlchart.v1.27").RefersTo = "=MPShuffleJD!`$E`$2:`$E`$34"
$wb.Names.Item("_xlchart.v1.28").RefersTo = "=MPShuffleJD!`$F`$1"
$wb.Names.Item("_xlchart.v1.29").RefersTo = "=MPShuffleJD!`$F`$2:`$F`$34"

# .30-.39: was SwapJD's block (auto-renamed to MPShuffleJD!...) -> reused for
# Top1Absent's block.
$wb.Names.Item("_xlchart.v1.30").RefersTo = "=Top1Absent!`$B`$1"
$wb.Names.Item("_xlchart.v1.31").RefersTo = "=Top1Absent!`$B`$2:`$B`$34"
$wb.Names.Item("_xlchart.v1.32").RefersTo = "=Top1Absent!`$C`$1"
$wb.Names.Item("_xlchart.v1.33").RefersTo = "=Top1Absent!`$C`$2:`$C`$34"
$wb.Names.Item("_xlchart.v1.34").RefersTo = "=Top1Absent!`$D`$1"
$wb.Names.Item("_xlchart.v1.35").RefersTo = "=Top1Absent!`$D`$2:`$D`$34"
$wb.Names.Item("_xlchart.v1.36").RefersTo = "=Top1Absent!`$E`$1"
$wb.Names.Item("_xlchart.v1.37").RefersTo = "=Top1Absent!`$E`$2:`$E`$34"
$wb.Names.Item("_xlchart.v1.38").RefersTo = "=Top1Absent!`$F`$1"
$wb.Names.Item("_xlchart.v1.39").RefersTo = "=Top1Absent!`$F`$2:`$F`$34"

# .40-.69: the now-redundant duplicate Top1Absent blocks -> drop them.
for ($i = 40; $i -le 69; $i++) {
    $wb.Names.Item("_xlchart.v1." + $i).Delete()
}

# Switch the active tab from Top1Absent to the renamed sheet.
$ws.Activate()
